$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B5 from 7.45 to 3.95
$ws.Range("B2").Value = 3.95
$ws.Range("B3").Value = 3.95
$ws.Range("B4").Value = 3.95
$ws.Range("B5").Value = 3.95

# Update B6:B7 from 1 to 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0

# Move the active cell selection to A7
$ws.Range("A7").Select()
